# DOM and Banner author ids added
# Updates publication rows (re-synced against refreshed OpenAlex/DOM-Banner
# author-id data) on Sheet1: rows 3-7 and 21-23 get corrected
# author/affiliation/work metadata. Each cell is temporarily forced to Text
# format before assignment (and the format is cleared right after) so that
# values which look like dates/numbers (publication_date, cited_by_count)
# are stored verbatim as text, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "William H. Frishman, Joseph S. Alpert"
$ws.Range("A3").ClearFormats()
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "New York Medical College/Westchester Medical Center, Valhalla; Supplements Editor,; The American Journal of Medicine; Department of Medicine, University of Arizona, Tucson; Editor in Chief,; The American Journal of Medicine"
$ws.Range("B3").ClearFormats()
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "https://openalex.org/W4318594540"
$ws.Range("C3").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "Commotio Cordis and the Triumph of Out-of-Hospital Cardiopulmonary Resuscitation"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2023-05-01"
$ws.Range("E3").ClearFormats()
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "https://doi.org/10.1016/j.amjmed.2023.01.007"
$ws.Range("H3").ClearFormats()
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = "3"
$ws.Range("M3").ClearFormats()
$ws.Range("O3").NumberFormat = "@"
$ws.Range("O3").Value = "https://pubmed.ncbi.nlm.nih.gov/36736646"
$ws.Range("O3").ClearFormats()
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = "https://doi.org/10.1016/j.amjmed.2023.01.007"
$ws.Range("P3").ClearFormats()

# Row 4
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "The American Journal of Medicine; University of Arizona School of Medicine, Tucson"
$ws.Range("B4").ClearFormats()
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "https://openalex.org/W4291377786"
$ws.Range("C4").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "Vasodilator Therapy in Hot Weather: A Warning"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2023-03-01"
$ws.Range("E4").ClearFormats()
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "https://doi.org/10.1016/j.amjmed.2022.08.003"
$ws.Range("H4").ClearFormats()
$ws.Range("M4").NumberFormat = "@"
$ws.Range("M4").Value = "2"
$ws.Range("M4").ClearFormats()
$ws.Range("O4").NumberFormat = "@"
$ws.Range("O4").Value = "https://pubmed.ncbi.nlm.nih.gov/35981649"
$ws.Range("O4").ClearFormats()
$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = "https://doi.org/10.1016/j.amjmed.2022.08.003"
$ws.Range("P4").ClearFormats()

# Row 5
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "Department of Medicine, University of Arizona, Tucson; Editor in Chief; The American Journal of Medicine"
$ws.Range("B5").ClearFormats()
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "https://openalex.org/W4310677416"
$ws.Range("C5").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "The Meaning of Life: To Serve Others"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2023-10-01"
$ws.Range("E5").ClearFormats()
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "https://doi.org/10.1016/j.amjmed.2022.11.006"
$ws.Range("H5").ClearFormats()
$ws.Range("O5").NumberFormat = "@"
$ws.Range("O5").Value = "https://pubmed.ncbi.nlm.nih.gov/36473501"
$ws.Range("O5").ClearFormats()
$ws.Range("P5").NumberFormat = "@"
$ws.Range("P5").Value = "https://doi.org/10.1016/j.amjmed.2022.11.006"
$ws.Range("P5").ClearFormats()

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "Joseph S. Alpert"
$ws.Range("A6").ClearFormats()
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "Editor in Chief,; The American Journal of Medicine; University of Arizona School of Medicine, Tucson"
$ws.Range("B6").ClearFormats()
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "https://openalex.org/W4321350851"
$ws.Range("C6").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "The Universal Panacea: Chicken Soup"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2023-12-01"
$ws.Range("E6").ClearFormats()
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "https://doi.org/10.1016/j.amjmed.2023.02.003"
$ws.Range("H6").ClearFormats()
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = "N/A"
$ws.Range("J6").ClearFormats()
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = "closed"
$ws.Range("K6").ClearFormats()
$ws.Range("O6").NumberFormat = "@"
$ws.Range("O6").Value = "https://pubmed.ncbi.nlm.nih.gov/36809812"
$ws.Range("O6").ClearFormats()
$ws.Range("P6").NumberFormat = "@"
$ws.Range("P6").Value = "https://doi.org/10.1016/j.amjmed.2023.02.003"
$ws.Range("P6").ClearFormats()

# Row 7
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "Editor in Chief,; The Amercian Journal of Medicine; University of Arizona School of Medicine, Tucson"
$ws.Range("B7").ClearFormats()
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://openalex.org/W4317930357"
$ws.Range("C7").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "Remarkable Advances in Clinical Medicine that Have Occurred Since I Was an Intern"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2023-06-01"
$ws.Range("E7").ClearFormats()
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "https://doi.org/10.1016/j.amjmed.2023.01.003"
$ws.Range("H7").ClearFormats()
$ws.Range("J7").NumberFormat = "@"
$ws.Range("J7").Value = "publishedVersion"
$ws.Range("J7").ClearFormats()
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "bronze"
$ws.Range("K7").ClearFormats()
$ws.Range("O7").NumberFormat = "@"
$ws.Range("O7").Value = "https://pubmed.ncbi.nlm.nih.gov/36707014"
$ws.Range("O7").ClearFormats()
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "https://doi.org/10.1016/j.amjmed.2023.01.003"
$ws.Range("P7").ClearFormats()

# Row 21
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "Susan Alpert, Kiril Solovey, Itzik Klein, Oren Salzman"
$ws.Range("A21").ClearFormats()
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "; ; ; "
$ws.Range("B21").ClearFormats()
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://openalex.org/W4386721834"
$ws.Range("C21").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "Inspection planning under execution uncertainty"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2023-09-12"
$ws.Range("E21").ClearFormats()
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "arXiv (Cornell University)"
$ws.Range("F21").ClearFormats()
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "Cornell University"
$ws.Range("G21").ClearFormats()
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "https://arxiv.org/abs/2309.06113"
$ws.Range("H21").ClearFormats()
$ws.Range("J21").NumberFormat = "@"
$ws.Range("J21").Value = "submittedVersion"
$ws.Range("J21").ClearFormats()
$ws.Range("K21").NumberFormat = "@"
$ws.Range("K21").Value = "green"
$ws.Range("K21").ClearFormats()
$ws.Range("O21").NumberFormat = "@"
$ws.Range("O21").Value = "NA"
$ws.Range("O21").ClearFormats()
$ws.Range("P21").NumberFormat = "@"
$ws.Range("P21").Value = "https://doi.org/10.48550/arxiv.2309.06113"
$ws.Range("P21").ClearFormats()

# Row 22
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "University of Arizona School of Medicine, Editor-in-Chief, The American Journal of Medicin, 1501 N. Campbell Avenue, Tucson, AZ 85724-5037"
$ws.Range("B22").ClearFormats()
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://openalex.org/W4387910389"
$ws.Range("C22").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "Twelve interesting biological tidbits"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2023-10-01"
$ws.Range("E22").ClearFormats()
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "https://doi.org/10.1016/j.amjmed.2023.10.015"
$ws.Range("H22").ClearFormats()
$ws.Range("J22").NumberFormat = "@"
$ws.Range("J22").Value = "publishedVersion"
$ws.Range("J22").ClearFormats()
$ws.Range("K22").NumberFormat = "@"
$ws.Range("K22").Value = "bronze"
$ws.Range("K22").ClearFormats()
$ws.Range("O22").NumberFormat = "@"
$ws.Range("O22").Value = "https://pubmed.ncbi.nlm.nih.gov/37879589"
$ws.Range("O22").ClearFormats()
$ws.Range("P22").NumberFormat = "@"
$ws.Range("P22").Value = "https://doi.org/10.1016/j.amjmed.2023.10.015"
$ws.Range("P22").ClearFormats()

# Row 23
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "Joseph S. Alpert"
$ws.Range("A23").ClearFormats()
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "University of Arizona School of Medicine, Tucson"
$ws.Range("B23").ClearFormats()
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://openalex.org/W4389849778"
$ws.Range("C23").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "Should Physicians and Nurses Ever Sit on the Patient's Bed?"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2023-12-01"
$ws.Range("E23").ClearFormats()
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "The American Journal of Medicine"
$ws.Range("F23").ClearFormats()
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "Elsevier BV"
$ws.Range("G23").ClearFormats()
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "https://doi.org/10.1016/j.amjmed.2023.12.013"
$ws.Range("H23").ClearFormats()
$ws.Range("J23").NumberFormat = "@"
$ws.Range("J23").Value = "publishedVersion"
$ws.Range("J23").ClearFormats()
$ws.Range("K23").NumberFormat = "@"
$ws.Range("K23").Value = "bronze"
$ws.Range("K23").ClearFormats()
$ws.Range("O23").NumberFormat = "@"
$ws.Range("O23").Value = "https://pubmed.ncbi.nlm.nih.gov/38110068"
$ws.Range("O23").ClearFormats()
$ws.Range("P23").NumberFormat = "@"
$ws.Range("P23").Value = "https://doi.org/10.1016/j.amjmed.2023.12.013"
$ws.Range("P23").ClearFormats()
